$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.051.38"
$ws.Range("E2").Value = "  +4.76%  "

$ws.Range("D3").Value = "2.245.73"
$ws.Range("E3").Value = "  +4.08%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.01"
$ws.Range("E5").Value = "  +3.99%  "

$ws.Range("E6").Value = "  +2.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.17"
$ws.Range("E7").Value = "  +8.91%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  +6.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.94"
$ws.Range("E10").Value = "  +6.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("E11").Value = "  +2.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.92"
$ws.Range("E12").Value = "  +4.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("E13").Value = "  +1.61%  "

$ws.Range("D14").Value = "2.583.65"
$ws.Range("E14").Value = "  +4.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.58"
$ws.Range("E15").Value = "  +2.87%  "

$ws.Range("D16").Value = "2.247.96"
$ws.Range("E16").Value = "  +5.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "42.948.05"
$ws.Range("E18").Value = "  +5.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000105"
$ws.Range("E19").Value = "  +5.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.04"
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").Value = "  +3.74%  "

$ws.Range("E22").Value = "  +7.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.71"
$ws.Range("E23").Value = "  +2.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  +16.66%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.86"
$ws.Range("E26").Value = "  +2.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.41"
$ws.Range("E27").Value = "  -0.59%  "

$ws.Range("E28").Value = "  +2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.70"
$ws.Range("E29").Value = "  +30.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.77"
$ws.Range("E30").Value = "  +2.19%  "

$ws.Range("E31").Value = "  -1.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.28"
$ws.Range("E32").Value = "  +2.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0794"
$ws.Range("E33").Value = "  +5.91%  "

$ws.Range("E34").Value = "  +4.55%  "

$ws.Range("E35").Value = "  +2.15%  "

$ws.Range("E36").Value = "  +7.89%  "

$ws.Range("E37").Value = "  +7.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0333"
$ws.Range("E38").Value = "  +20.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.02"
$ws.Range("E39").Value = "  +14.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.12"
$ws.Range("E40").Value = "  +3.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.47"
$ws.Range("E41").Value = "  +3.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.205"
$ws.Range("E42").Value = "  +10.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.26"
$ws.Range("E43").Value = "  +2.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.77"
$ws.Range("E44").Value = "  +8.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0989"
$ws.Range("E47").Value = "  +3.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +11.32%  "

$ws.Range("E49").Value = "  +4.01%  "

$ws.Range("E50").Value = "  +4.31%  "

$ws.Range("D51").Value = "2.457.03"
$ws.Range("E51").Value = "  +4.21%  "

# Row 45: now WOONetwork (was FraxShare)
$ws.Range("B45").Value = "WOONetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.489"
$ws.Range("E45").Value = "  +33.78%  "

# Row 46: now FraxShare (was WOONetwork)
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.69"
$ws.Range("E46").Value = "  +6.05%  "
